$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 401.43332
$ws.Cells.Item(28, 9).Value = 374.45456
$ws.Cells.Item(28, 11).Value = 374.45456
$ws.Cells.Item(28, 13).Value = 110.54544
$ws.Cells.Item(32, 8).Value = 2948.5715
$ws.Cells.Item(32, 9).Value = 3333
$ws.Cells.Item(32, 10).Value = 2884.5
$ws.Cells.Item(32, 11).Value = 3333
$ws.Cells.Item(32, 12).Value = 2884.5
$ws.Cells.Item(32, 13).Value = -3007
$ws.Cells.Item(32, 14).Value = -3536.5
$ws.Cells.Item(62, 8).Value = 10633.77
$ws.Cells.Item(62, 10).Value = 5552.4
$ws.Cells.Item(62, 12).Value = 5552.4
$ws.Cells.Item(62, 14).Value = -6800.4
$ws.Cells.Item(65, 8).Value = 10633.77
$ws.Cells.Item(65, 10).Value = 5552.4
$ws.Cells.Item(65, 12).Value = 27762
$ws.Cells.Item(65, 14).Value = -34002
$ws.Cells.Item(86, 8).Value = 3092.8333
$ws.Cells.Item(86, 9).Value = 2689.5
$ws.Cells.Item(86, 11).Value = 2689.5
$ws.Cells.Item(86, 13).Value = -1566.5
$ws.Cells.Item(89, 8).Value = 3092.8333
$ws.Cells.Item(89, 9).Value = 2689.5
$ws.Cells.Item(89, 11).Value = 13447.5
$ws.Cells.Item(89, 13).Value = -7831.5
$ws.Cells.Item(98, 8).Value = 3302.1428
$ws.Cells.Item(98, 9).Value = 1754.7142
$ws.Cells.Item(98, 10).Value = 4849.5713
$ws.Cells.Item(98, 11).Value = 1754.7142
$ws.Cells.Item(98, 12).Value = 4849.5713
$ws.Cells.Item(98, 13).Value = -256.7141999999999
$ws.Cells.Item(98, 14).Value = -7845.5713
$ws.Cells.Item(100, 8).Value = 2832.8206
$ws.Cells.Item(100, 9).Value = 2351.4517
$ws.Cells.Item(100, 10).Value = 4698.125
$ws.Cells.Item(100, 11).Value = 2351.4517
$ws.Cells.Item(100, 12).Value = 4698.125
$ws.Cells.Item(100, 13).Value = -1810.4517
$ws.Cells.Item(100, 14).Value = -5780.125
$ws.Cells.Item(111, 8).Value = 2349.625
$ws.Cells.Item(111, 9).Value = 2324.5454
$ws.Cells.Item(111, 10).Value = 2404.8
$ws.Cells.Item(111, 11).Value = 6973.6362
$ws.Cells.Item(111, 12).Value = 7214.400000000001
$ws.Cells.Item(111, 13).Value = -3906.6362
$ws.Cells.Item(111, 14).Value = -13348.4
$ws.Cells.Item(122, 8).Value = 3302.1428
$ws.Cells.Item(122, 9).Value = 1754.7142
$ws.Cells.Item(122, 10).Value = 4849.5713
$ws.Cells.Item(122, 11).Value = 5264.142599999999
$ws.Cells.Item(122, 12).Value = 14548.7139
$ws.Cells.Item(122, 13).Value = -2814.142599999999
$ws.Cells.Item(122, 14).Value = -19448.7139
$ws.Cells.Item(132, 8).Value = 5538.1025
$ws.Cells.Item(132, 9).Value = 4345.6284
$ws.Cells.Item(132, 10).Value = 15972.25
$ws.Cells.Item(132, 11).Value = 13036.8852
$ws.Cells.Item(132, 12).Value = 47916.75
$ws.Cells.Item(132, 13).Value = -10506.8852
$ws.Cells.Item(132, 14).Value = -52976.75
$ws.Cells.Item(138, 8).Value = 5455.6
$ws.Cells.Item(138, 9).Value = 11192.667
$ws.Cells.Item(138, 10).Value = 4673.273
$ws.Cells.Item(138, 11).Value = 33578.001
$ws.Cells.Item(138, 12).Value = 14019.819
$ws.Cells.Item(138, 13).Value = -28438.001
$ws.Cells.Item(138, 14).Value = -24299.819
$ws.Cells.Item(141, 8).Value = 2690
$ws.Cells.Item(141, 9).Value = 2745
$ws.Cells.Item(141, 10).Value = 2140
$ws.Cells.Item(141, 11).Value = 8235
$ws.Cells.Item(141, 12).Value = 6420
$ws.Cells.Item(141, 13).Value = -3055
$ws.Cells.Item(141, 14).Value = -16780

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4298.9165
$ws.Cells.Item(45, 9).Value = 3348.7144
$ws.Cells.Item(45, 10).Value = 5629.2
$ws.Cells.Item(45, 11).Value = 3348.7144
$ws.Cells.Item(45, 12).Value = 5629.2
$ws.Cells.Item(45, 13).Value = -2971.7144
$ws.Cells.Item(45, 14).Value = -6383.2
$ws.Cells.Item(74, 8).Value = 2552.9167
$ws.Cells.Item(74, 9).Value = 2421.3635
$ws.Cells.Item(74, 10).Value = 4000
$ws.Cells.Item(74, 11).Value = 2421.3635
$ws.Cells.Item(74, 12).Value = 4000
$ws.Cells.Item(74, 13).Value = -1547.3635
$ws.Cells.Item(74, 14).Value = -5748
$ws.Cells.Item(77, 8).Value = 2552.9167
$ws.Cells.Item(77, 9).Value = 2421.3635
$ws.Cells.Item(77, 10).Value = 4000
$ws.Cells.Item(77, 11).Value = 12106.8175
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = -7738.817499999999
$ws.Cells.Item(77, 14).Value = -28736
$ws.Cells.Item(102, 8).Value = 5566.3335
$ws.Cells.Item(102, 9).Value = 4676
$ws.Cells.Item(102, 11).Value = 4676
$ws.Cells.Item(102, 13).Value = -3054
$ws.Cells.Item(110, 8).Value = 3845.2632
$ws.Cells.Item(110, 9).Value = 3504.6428
$ws.Cells.Item(110, 10).Value = 4799
$ws.Cells.Item(110, 11).Value = 3504.6428
$ws.Cells.Item(110, 12).Value = 4799
$ws.Cells.Item(110, 13).Value = -1459.6428
$ws.Cells.Item(110, 14).Value = -8889

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 46611
$ws.Cells.Item(20, 9).Value = 55267.58
$ws.Cells.Item(20, 11).Value = 55267.58
$ws.Cells.Item(20, 13).Value = -55020.58
$ws.Cells.Item(105, 8).Value = 3745
$ws.Cells.Item(105, 9).Value = 3745
$ws.Cells.Item(105, 11).Value = 3745
$ws.Cells.Item(105, 13).Value = -1998
$ws.Cells.Item(107, 8).Value = 2348.8333
$ws.Cells.Item(107, 9).Value = 1668.8889
$ws.Cells.Item(107, 10).Value = 4388.6665
$ws.Cells.Item(107, 11).Value = 1668.8889
$ws.Cells.Item(107, 12).Value = 4388.6665
$ws.Cells.Item(107, 13).Value = 251.1111000000001
$ws.Cells.Item(107, 14).Value = -8228.666499999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1679.0769
$ws.Cells.Item(31, 9).Value = 1552
$ws.Cells.Item(31, 10).Value = 2653.3333
$ws.Cells.Item(31, 11).Value = 1552
$ws.Cells.Item(31, 12).Value = 2653.3333
$ws.Cells.Item(31, 13).Value = -1257
$ws.Cells.Item(31, 14).Value = -3243.3333
$ws.Cells.Item(34, 8).Value = 1679.0769
$ws.Cells.Item(34, 9).Value = 1552
$ws.Cells.Item(34, 10).Value = 2653.3333
$ws.Cells.Item(34, 11).Value = 1552
$ws.Cells.Item(34, 12).Value = 2653.3333
$ws.Cells.Item(34, 13).Value = -1350
$ws.Cells.Item(34, 14).Value = -3057.3333
$ws.Cells.Item(42, 8).Value = 3330.6
$ws.Cells.Item(42, 9).Value = 3330.6
$ws.Cells.Item(42, 11).Value = 3330.6
$ws.Cells.Item(42, 13).Value = -2737.6
$ws.Cells.Item(44, 8).Value = 4999
$ws.Cells.Item(44, 9).Value = 4999
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 4999
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).ClearContents()
$ws.Cells.Item(44, 14).Value = -4557
$ws.Cells.Item(62, 8).Value = 3462.5
$ws.Cells.Item(62, 9).Value = 3194
$ws.Cells.Item(62, 10).Value = 3999.5
$ws.Cells.Item(62, 11).Value = 3194
$ws.Cells.Item(62, 12).Value = 3999.5
$ws.Cells.Item(62, 13).Value = -2570
$ws.Cells.Item(62, 14).Value = -5247.5
$ws.Cells.Item(65, 8).Value = 3462.5
$ws.Cells.Item(65, 9).Value = 3194
$ws.Cells.Item(65, 10).Value = 3999.5
$ws.Cells.Item(65, 11).Value = 15970
$ws.Cells.Item(65, 12).Value = 19997.5
$ws.Cells.Item(65, 13).Value = -12850
$ws.Cells.Item(65, 14).Value = -26237.5
$ws.Cells.Item(94, 8).Value = 1729.1765
$ws.Cells.Item(94, 9).Value = 1728.4286
$ws.Cells.Item(94, 10).Value = 1729.7
$ws.Cells.Item(94, 11).Value = 1728.4286
$ws.Cells.Item(94, 12).Value = 1729.7
$ws.Cells.Item(94, 13).Value = -1277.4286
$ws.Cells.Item(94, 14).Value = -2631.7
$ws.Cells.Item(132, 8).Value = 1832.5454
$ws.Cells.Item(132, 9).Value = 1832.5454
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 5497.6362
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -2967.6362
$ws.Cells.Item(134, 8).Value = 58865
$ws.Cells.Item(134, 9).Value = 61974.824
$ws.Cells.Item(134, 11).Value = 185924.472
$ws.Cells.Item(134, 13).Value = -183389.472

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 9610.888999999999
$ws.Cells.Item(19, 9).Value = 8250
$ws.Cells.Item(19, 10).Value = 9999.714
$ws.Cells.Item(19, 11).Value = 24750
$ws.Cells.Item(19, 12).Value = 29999.142
$ws.Cells.Item(19, 13).Value = -24576
$ws.Cells.Item(19, 14).Value = -30347.142
$ws.Cells.Item(33, 8).Value = 642.7143
$ws.Cells.Item(33, 10).Value = 1066.3334
$ws.Cells.Item(33, 12).Value = 6398.0004
$ws.Cells.Item(33, 14).Value = -6964.0004
$ws.Cells.Item(37, 8).Value = 60011.11
$ws.Cells.Item(37, 10).Value = 60011.11
$ws.Cells.Item(37, 12).Value = 180033.33
$ws.Cells.Item(37, 14).Value = -180257.33
$ws.Cells.Item(39, 8).Value = 4739.1875
$ws.Cells.Item(39, 10).Value = 4739.1875
$ws.Cells.Item(39, 12).Value = 14217.5625
$ws.Cells.Item(39, 14).Value = -14805.5625
$ws.Cells.Item(40, 8).Value = 52.555557
$ws.Cells.Item(40, 9).Value = 24.714285
$ws.Cells.Item(40, 10).Value = 150
$ws.Cells.Item(40, 11).Value = 98.85714
$ws.Cells.Item(40, 12).Value = 600
$ws.Cells.Item(40, 13).Value = -29.85714
$ws.Cells.Item(40, 14).Value = -738
$ws.Cells.Item(55, 8).Value = 4165.222
$ws.Cells.Item(55, 10).Value = 6064.8335
$ws.Cells.Item(55, 12).Value = 18194.5005
$ws.Cells.Item(55, 14).Value = -18548.5005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4964.909
$ws.Cells.Item(70, 9).Value = 4771.1665
$ws.Cells.Item(70, 11).Value = 4771.1665
$ws.Cells.Item(70, 13).Value = -4501.1665
$ws.Cells.Item(73, 8).Value = 4964.909
$ws.Cells.Item(73, 9).Value = 4771.1665
$ws.Cells.Item(73, 11).Value = 4771.1665
$ws.Cells.Item(73, 13).Value = -3835.1665

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2119.1191
$ws.Cells.Item(16, 9).Value = 1735.8823
$ws.Cells.Item(16, 11).Value = 1735.8823
$ws.Cells.Item(16, 13).Value = -1565.8823
$ws.Cells.Item(46, 8).Value = 2809.625
$ws.Cells.Item(46, 9).Value = 962
$ws.Cells.Item(46, 10).Value = 3179.15
$ws.Cells.Item(46, 11).Value = 962
$ws.Cells.Item(46, 12).Value = 3179.15
$ws.Cells.Item(46, 13).Value = -774
$ws.Cells.Item(46, 14).Value = -3555.15
$ws.Cells.Item(127, 8).Value = 48098.6
$ws.Cells.Item(127, 10).Value = 48098.6
$ws.Cells.Item(127, 12).Value = 48098.6
$ws.Cells.Item(127, 14).Value = -58018.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 99313.27
$ws.Cells.Item(62, 9).Value = 8999
$ws.Cells.Item(62, 10).Value = 119383.11
$ws.Cells.Item(62, 11).Value = 8999
$ws.Cells.Item(62, 12).Value = 119383.11
$ws.Cells.Item(62, 13).Value = -8375
$ws.Cells.Item(62, 14).Value = -120631.11
$ws.Cells.Item(65, 8).Value = 99313.27
$ws.Cells.Item(65, 9).Value = 8999
$ws.Cells.Item(65, 10).Value = 119383.11
$ws.Cells.Item(65, 11).Value = 44995
$ws.Cells.Item(65, 12).Value = 596915.55
$ws.Cells.Item(65, 13).Value = -41875
$ws.Cells.Item(65, 14).Value = -603155.55
$ws.Cells.Item(113, 8).Value = 4588
$ws.Cells.Item(113, 9).Value = 2412.6
$ws.Cells.Item(113, 11).Value = 7237.799999999999
$ws.Cells.Item(113, 13).Value = -5067.799999999999
$ws.Cells.Item(122, 8).Value = 7232.9546
$ws.Cells.Item(122, 9).Value = 8091.263
$ws.Cells.Item(122, 11).Value = 24273.789
$ws.Cells.Item(122, 13).Value = -21823.789
$ws.Cells.Item(126, 8).Value = 73439
$ws.Cells.Item(126, 9).Value = 89719.914
$ws.Cells.Item(126, 11).Value = 269159.742
$ws.Cells.Item(126, 13).Value = -266689.742
$ws.Cells.Item(132, 8).Value = 52205.65
$ws.Cells.Item(132, 9).Value = 57634.61
$ws.Cells.Item(132, 11).Value = 172903.83
$ws.Cells.Item(132, 13).Value = -170373.83
